$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row before row 2 (shifts existing rows 2..22 down to 3..23)
$ws.Rows.Item(2).Insert()

# Copy formatting from row 3 (the row that used to be row 2, now shifted down)
# down into the newly inserted blank row 2 so styles match the rest of the table.
$ws.Range("A3:F3").Copy()
$ws.Range("A2:F2").PasteSpecial(-4122)

# Fill in the new row 2 data
$ws.Range("A2").Value = 22
$ws.Range("B2").Value = "2. P0610 (99.85% min) /P1020/ EC Grade Ingot & Sow 99.7% (min) / Cast Bar"
$ws.Range("C2").Value = "P1020"
$ws.Range("D2").Value = 267.25
$ws.Range("E2").Value = "13.08.2025"
$ws.Range("F2").Value = "https://www.hindalco.com/Upload/PDF/primary-ready-reckoner-13-august-2025.pdf"

# Fix up hyperlinks: remove the old one (which did not shift with the row insert)
# and re-add hyperlinks at the correct cells: F2 (new) and F3 (shifted from old F2).
$ws.Hyperlinks.Delete()
$ws.Hyperlinks.Add($ws.Range("F2"), "https://www.hindalco.com/Upload/PDF/primary-ready-reckoner-13-august-2025.pdf")
$ws.Hyperlinks.Add($ws.Range("F3"), "https://www.hindalco.com/Upload/PDF/primary-ready-reckoner-12-august-2025.pdf")

# Adding a hyperlink auto-applies Excel's built-in "Hyperlink" style (blue/underline),
# but the original workbook keeps the plain table style on these cells, so restore it.
$ws.Range("A2").Copy()
$ws.Range("F2").PasteSpecial(-4122)
$ws.Range("A3").Copy()
$ws.Range("F3").PasteSpecial(-4122)
